$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.896310567855835
$ws.Range("B1").Value = 1.691284537315369
$ws.Range("C1").Value = 4.235966205596924
$ws.Range("D1").Value = 3.325758934020996
$ws.Range("E1").Value = 0.43682861328125
